$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A8").Value = "ShonJob1"
$ws.Range("B8").Value = "Test1!!!"
$ws.Range("A9").Value = "Habuf"
$ws.Range("B9").Value = "Habuf12!"
